# Applies the tracked changes to s3_good_trials_GraspObject_Shuffled.xlsx:
#  - mark rows 4 and 6 (columns A) as "not included"
#  - add a new data row (row 17) for session 20240829
#  - shrink column widths B:G
#  - move the active selection to E19

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Mark the two sessions in column A as "not included" ---
$ws.Cells.Item(4, 1).Value = "not included"
$ws.Cells.Item(6, 1).Value = "not included"

# --- Append the new session row (row 17) ---
$ws.Cells.Item(17, 1).Value = 20240829
$ws.Cells.Item(17, 2).Value = 3
$ws.Cells.Item(17, 3).Value = 5
$ws.Cells.Item(17, 4).Value = 6
$ws.Cells.Item(17, 5).Value = 8
$ws.Cells.Item(17, 6).Value = 4
$ws.Cells.Item(17, 7).Value = 7

# --- Narrow columns B:G to their new widths ---
$ws.Columns.Item(2).ColumnWidth = 21.1666666666667
$ws.Columns.Item(3).ColumnWidth = 17.5
$ws.Columns.Item(4).ColumnWidth = 20
$ws.Columns.Item(5).ColumnWidth = 21.3333333333333
$ws.Columns.Item(6).ColumnWidth = 22.3333333333333
$ws.Columns.Item(7).ColumnWidth = 18.5

# --- Move the selection like the author left it ---
$ws.Range("E19").Select() | Out-Null
